$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: new "Parameters set" label for the STAR rerun / novoalign block
$ws.Range("A17").Value = "Parameters set: --outFilterScoreMinOverLread 0 --outFilterMatchNminOverLread 0 --outFilterMatchNmin 0 --outFilterMismatchNmax 2"

# Row 18: header row (same labels as row 3)
$ws.Range("A18").Value = "Sample"
$ws.Range("B18").Value = "# of input reads"
$ws.Range("C18").Value = "Average input read length"
$ws.Range("D18").Value = "% uniquely mapped reads"
$ws.Range("E18").Value = "Average insert length"
$ws.Range("F18").Value = "% multimapping reads"
$ws.Range("G18").Value = "% unmapped reads (mismatch)"
$ws.Range("H18").Value = "% unmapped reads (too short)"
$ws.Range("I18").Value = "% unmapped reads (other)"

# Data rows 19-30
$data = @(
  @{ Row=19; A="SRP179837_leaf1"; B=22595120; C=300; D=0.79969999999999997; E=2.27;               F=0.1973;               G=0; H=0; I=0.00040000000000000002 },
  @{ Row=20; A="SRP179837_leaf2"; B=25589730; C=300; D=0.80400000000000005; E=2.44;               F=0.19259999999999999;  G=0; H=0; I=0.00050000000000000001 },
  @{ Row=21; A="SRP179837_leaf3"; B=48528038; C=300; D=0.7268;              E=2.27;               F=0.26469999999999999;  G=0; H=0; I=0.00040000000000000002 },
  @{ Row=22; A="SRP179837_root1"; B=41235603; C=300; D=0.76500000000000001; E=2.2799999999999998; F=0.23150000000000001;  G=0; H=0; I=0.00059999999999999995 },
  @{ Row=23; A="SRP179837_root2"; B=20852584; C=300; D=0.78710000000000002; E=2.4500000000000002; F=0.2094;               G=0; H=0; I=0.00050000000000000001 },
  @{ Row=24; A="SRP179837_root3"; B=18571158; C=300; D=0.79220000000000002; E=2.2799999999999998; F=0.20449999999999999;  G=0; H=0; I=0.00029999999999999997 },
  @{ Row=25; A="RAC_leaf1";       B=38680165; C=202; D=0.69489999999999996; E=2.76;               F=0.29849999999999999;  G=0; H=0; I=0.0001 },
  @{ Row=26; A="RAC_leaf2";       B=39743330; C=202; D=0.69;                E=2.69;               F=0.30309999999999998;  G=0; H=0; I=0.0001 },
  @{ Row=27; A="RAC_leaf3";       B=31208809; C=202; D=0.71150000000000002; E=2.6;                F=0.28270000000000001;  G=0; H=0; I=0.0001 },
  @{ Row=28; A="RAC_root1";       B=32802143; C=202; D=0.70830000000000004; E=2.5299999999999998; F=0.28549999999999998;  G=0; H=0; I=0.00020000000000000001 },
  @{ Row=29; A="RAC_root2";       B=29613922; C=202; D=0.73199999999999998; E=2.61;               F=0.2641;               G=0; H=0; I=0.00020000000000000001 },
  @{ Row=30; A="RAC_root3";       B=35774916; C=202; D=0.73040000000000005; E=2.67;               F=0.2651;               G=0; H=0; I=0.00020000000000000001 }
)

foreach ($rec in $data) {
  $r = $rec.Row
  $ws.Range("A$r").Value = $rec.A
  $ws.Range("B$r").Value = $rec.B
  $ws.Range("C$r").Value = $rec.C
  $ws.Range("D$r").Value = $rec.D
  $ws.Range("E$r").Value = $rec.E
  $ws.Range("F$r").Value = $rec.F
  $ws.Range("G$r").Value = $rec.G
  $ws.Range("H$r").Value = $rec.H
  $ws.Range("I$r").Value = $rec.I

  $ws.Range("D$r").NumberFormat = "0.00%"
  $ws.Range("F$r").NumberFormat = "0.00%"
  $ws.Range("G$r").NumberFormat = "0%"
  $ws.Range("H$r").NumberFormat = "0%"
  $ws.Range("I$r").NumberFormat = "0.00%"
}

# Update selection to match the committed view state
$ws.Range("C16").Select()
